$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("IFSRF")

# --- Row 11 (hydrogen if): the old "hydrogen if" recipient flags move to the
# new "green hydrogen if" row, so zero them out here first.
$ws.Range("C11").Value = 0
$ws.Range("D11").Value = 0
$ws.Range("K11").Value = 0

# --- New header columns L (green hydrogen if) and M (low carbon hydrogen if)
$ws.Range("L1").Value = "green hydrogen if"
$ws.Range("M1").Value = "low carbon hydrogen if"
$ws.Range("K1").Copy() | Out-Null
$ws.Range("L1:M1").PasteSpecial(-4122) | Out-Null  # xlPasteFormats

# --- Fill new L,M columns for rows 2-11 with 0, copying number column style
$ws.Range("K2:K11").Copy() | Out-Null
$ws.Range("L2:M11").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("L2:M11").Value = 0

# --- New row 12: "green hydrogen if" (absorbs the old hydrogen-if flags and
# picks up new green-/low-carbon-hydrogen "to" flags as well)
$ws.Range("A11:K11").Copy() | Out-Null
$ws.Range("A12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("L11:M11").Copy() | Out-Null
$ws.Range("L12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A12").Value = "green hydrogen if"
$ws.Range("B12:M12").Value = 0
$ws.Range("C12").Value = 1
$ws.Range("D12").Value = 1
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 1
$ws.Range("M12").Value = 1

# --- New row 13: "low carbon hydrogen if"
$ws.Range("A12:M12").Copy() | Out-Null
$ws.Range("A13").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("A13").Value = "low carbon hydrogen if"
$ws.Range("B13:M13").Value = 0

# --- Row 14: formatting-only remnants in L14/M14 (no values), matching the
# trailing empty-styled row left behind by the fill-down.
$ws.Range("L13:M13").Copy() | Out-Null
$ws.Range("L14").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$ws.Range("L14:M14").ClearContents() | Out-Null

# --- Sheet view / selection bookkeeping: the author ended editing on the
# About sheet (now the active tab) with D12 left selected on IFSRF behind it.
$ws.Range("D12").Select() | Out-Null
$wb.Worksheets.Item("About").Activate() | Out-Null
$wb.Worksheets.Item("About").Range("E12").Select() | Out-Null
